$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.981.45"
$ws.Range("E2").Value = "  -3.31%  "

# Row 3
$ws.Range("D3").Value = "3.354.69"
$ws.Range("E3").Value = "  -3.75%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "'558.11"
$ws.Range("E5").Value = "  +0.40%  "

# Row 6
$ws.Range("D6").Value = "'172.48"
$ws.Range("E6").Value = "  -6.46%  "

# Row 7
$ws.Range("D7").Value = "'0.613"
$ws.Range("E7").Value = "  -5.82%  "

# Row 8
$ws.Range("D8").Value = "3.338.03"
$ws.Range("E8").Value = "  -4.08%  "

# Row 9
$ws.Range("E9").Value = "  +0.09%  "

# Row 10
$ws.Range("D10").Value = "'0.619"
$ws.Range("E10").Value = "  -2.38%  "

# Row 11
$ws.Range("D11").Value = "'0.151"
$ws.Range("E11").Value = "  -1.23%  "

# Row 12
$ws.Range("D12").Value = "'54.09"
$ws.Range("E12").Value = "  -0.50%  "

# Row 13
$ws.Range("D13").Value = "'0.0000266"
$ws.Range("E13").Value = "  -1.32%  "

# Row 14
$ws.Range("D14").Value = "'8.92"
$ws.Range("E14").Value = "  -3.70%  "

# Row 15
$ws.Range("D15").Value = "3.882.06"
$ws.Range("E15").Value = "  -4.25%  "

# Row 16
$ws.Range("D16").Value = "'0.118"
$ws.Range("E16").Value = "  -2.85%  "

# Row 17
$ws.Range("D17").Value = "3.355.72"
$ws.Range("E17").Value = "  -3.98%  "

# Row 18
$ws.Range("D18").Value = "'17.72"
$ws.Range("E18").Value = "  -4.21%  "

# Row 19
$ws.Range("D19").Value = "'11.70"
$ws.Range("E19").Value = "  -2.29%  "

# Row 20
$ws.Range("D20").Value = "63.961.41"
$ws.Range("E20").Value = "  -3.44%  "

# Row 21
$ws.Range("D21").Value = "'0.974"
$ws.Range("E21").Value = "  -1.58%  "

# Row 22
$ws.Range("D22").Value = "'404.31"
$ws.Range("E22").Value = "  -4.01%  "

# Row 23
$ws.Range("D23").Value = "'4.10"
$ws.Range("E23").Value = "  +1.36%  "

# Row 24
$ws.Range("D24").Value = "'4.30"
$ws.Range("E24").Value = "  +4.38%  "

# Row 25
$ws.Range("D25").Value = "'82.82"
$ws.Range("E25").Value = "  -4.09%  "

# Row 26
$ws.Range("D26").Value = "'13.21"
$ws.Range("E26").Value = "  +7.75%  "

# Row 27
$ws.Range("D27").Value = "'10.76"
$ws.Range("E27").Value = "  -1.40%  "

# Row 28
$ws.Range("D28").Value = "'2.75"
$ws.Range("E28").Value = "  -4.96%  "

# Row 29
$ws.Range("D29").Value = "'8.78"
$ws.Range("E29").Value = "  -3.18%  "

# Row 30
$ws.Range("D30").Value = "'29.23"
$ws.Range("E30").Value = "  -2.84%  "

# Row 31
$ws.Range("D31").Value = "'6.48"
$ws.Range("E31").Value = "  -1.41%  "

# Row 32
$ws.Range("D32").Value = "'586.26"
$ws.Range("E32").Value = "  -7.02%  "

# Row 33
$ws.Range("D33").Value = "'11.36"
$ws.Range("E33").Value = "  -2.85%  "

# Row 34
$ws.Range("D34").Value = "'0.107"
$ws.Range("E34").Value = "  -3.23%  "

# Row 35
$ws.Range("D35").Value = "'57.93"
$ws.Range("E35").Value = "  -3.21%  "

# Row 36
$ws.Range("E36").Value = "  +1.30%  "

# Row 37
$ws.Range("E37").Value = "  +0.11%  "

# Row 38
$ws.Range("D38").Value = "'35.77"
$ws.Range("E38").Value = "  -4.90%  "

# Row 39
$ws.Range("D39").Value = "'3.43"
$ws.Range("E39").Value = "  +0.79%  "

# Row 40
$ws.Range("D40").Value = "0.0₃0748"
$ws.Range("E40").Value = "  -7.20%  "

# Row 41
$ws.Range("D41").Value = "'0.370"
$ws.Range("E41").Value = "  -3.61%  "

# Row 42
$ws.Range("D42").Value = "3.151.67"
$ws.Range("E42").Value = "  +1.44%  "

# Row 43
$ws.Range("E43").Value = "  -0.15%  "

# Row 44
$ws.Range("D44").Value = "'2.85"
$ws.Range("E44").Value = "  +0.33%  "

# Row 45
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "'3.24"
$ws.Range("E45").Value = "  -3.20%  "

# Row 46
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "'2.49"
$ws.Range("E46").Value = "  -3.86%  "

# Row 47
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0407"
$ws.Range("E47").Value = "  -1.53%  "

# Row 48
$ws.Range("D48").Value = "'2.61"
$ws.Range("E48").Value = "  -3.93%  "

# Row 49
$ws.Range("D49").Value = "'0.129"
$ws.Range("E49").Value = "  -4.48%  "

# Row 50
$ws.Range("D50").Value = "'133.23"
$ws.Range("E50").Value = "  -4.11%  "

# Row 51
$ws.Range("D51").Value = "'8.11"
$ws.Range("E51").Value = "  -4.31%  "
